# Update gh-pages to output generated at 456a3b4
# Apply numeric updates to the "想去人数" (interest count) column F
# on the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value  = 170
$ws1.Range("F8").Value  = 6337
$ws1.Range("F9").Value  = 73
$ws1.Range("F12").Value = 5425
$ws1.Range("F24").Value = 3888
$ws1.Range("F25").Value = 168

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value  = 170
$ws4.Range("F9").Value  = 6337
$ws4.Range("F10").Value = 73
$ws4.Range("F13").Value = 5425
$ws4.Range("F25").Value = 3888
$ws4.Range("F27").Value = 168
